# Add new test case "verifyThatUserCannotChangePasswordWhenTheNewPasswordsDoNotMatch"
# under UserAccountManagement: one row in RUNMANAGER (sheet1) plus two data rows
# (chrome/firefox) in USERACCOUNTMANAGEMENTDATA (sheet3).

$wb = $excel.ActiveWorkbook

# --- Sheet "RUNMANAGER": register the new test case (row 8) ---
# Write the test name first so it claims shared-string slot 38 before the
# USERACCOUNTMANAGEMENTDATA sheet allocates "admin132".
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Cells.Item(8, 1).Value = "verifyThatUserCannotChangePasswordWhenTheNewPasswordsDoNotMatch"
$ws1.Cells.Item(8, 2).Value = "To check this test is executed"
$ws1.Cells.Item(8, 3).Value = "yes"

# --- Sheet "USERACCOUNTMANAGEMENTDATA": two data-driven rows (chrome + firefox) ---
$ws3 = $wb.Worksheets.Item("USERACCOUNTMANAGEMENTDATA")

$ws3.Cells.Item(6, 1).Value = "verifyThatUserCannotChangePasswordWhenTheNewPasswordsDoNotMatch"
$ws3.Cells.Item(6, 2).Value = "yes"
$ws3.Cells.Item(6, 3).Value = "Admin"
$ws3.Cells.Item(6, 4).Value = "admin123"
$ws3.Cells.Item(6, 5).Value = "Sunil"
$ws3.Cells.Item(6, 6).Value = "chrome"
$ws3.Cells.Item(6, 7).Value = "admin123"
$ws3.Cells.Item(6, 8).Value = "admin111"
$ws3.Cells.Item(6, 9).Value = "admin132"

$ws3.Cells.Item(7, 1).Value = "verifyThatUserCannotChangePasswordWhenTheNewPasswordsDoNotMatch"
$ws3.Cells.Item(7, 2).Value = "yes"
$ws3.Cells.Item(7, 3).Value = "Admin"
$ws3.Cells.Item(7, 4).Value = "admin123"
$ws3.Cells.Item(7, 5).Value = "Sunil"
$ws3.Cells.Item(7, 6).Value = "firefox"
$ws3.Cells.Item(7, 7).Value = "admin123"
$ws3.Cells.Item(7, 8).Value = "admin111"
$ws3.Cells.Item(7, 9).Value = "admin132"

# --- back on RUNMANAGER: priority "7" / count "1" stored as text (quote-prefixed,
#     matching how every other value in these columns is stored) ---
$ws1.Cells.Item(8, 4).Value = "'7"
$ws1.Cells.Item(8, 5).Value = "'1"

# Column A on RUNMANAGER needs to widen to fit the new (longer) test name.
$ws1.Columns.Item(1).ColumnWidth = 60.83

# Restore each sheet's own selection (without leaving USERACCOUNTMANAGEMENTDATA
# as the active tab).
$ws3.Range("A7").Select()
$ws1.Range("C9").Select()

Write-Output "Added verifyThatUserCannotChangePasswordWhenTheNewPasswordsDoNotMatch test rows."
